{"js": "// Apply the four text replacements described by the diff.\nconst replacements = [\n  {\n    find: \"Total time calculated in dataset: 30 days 23:45:00\",\n    replace: \"Total time in days calculated in dataset: 30.99\"\n  },\n  {\n    find: \"Total time in hours for when fault flag is True: 10.0\",\n    replace: \"Total time in hours for when fault flag is True: 274.0\"\n  },\n  {\n    find: \"Fan is appears to generate good duct static pressure (GOOD)\",\n    replace: \"The percent True metric that represents the amount of time for when the fault flag is True is high indicating the fan is running at high speeds and appearing to not generate good duct static pressure\"\n  },\n  {\n    find: \"Report generated: Thu Dec 22 08:11:42 2022\",\n    replace: \"Report generated: Sun Dec 25 08:39:48 2022\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the four text replacements described by the diff using Find/Replace.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0  # wdFindStop\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute(\n        $findText,       # FindText\n        $true,           # MatchCase\n        $false,          # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        1,               # Wrap (wdFindContinue)\n        $false,          # Format\n        $replaceText,    # ReplaceWith\n        2                # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\nReplace-Text \"Total time calculated in dataset: 30 days 23:45:00\" \"Total time in days calculated in dataset: 30.99\"\nReplace-Text \"Total time in hours for when fault flag is True: 10.0\" \"Total time in hours for when fault flag is True: 274.0\"\nReplace-Text \"Fan is appears to generate good duct static pressure (GOOD)\" \"The percent True metric that represents the amount of time for when the fault flag is True is high indicating the fan is running at high speeds and appearing to not generate good duct static pressure\"\nReplace-Text \"Report generated: Thu Dec 22 08:11:42 2022\" \"Report generated: Sun Dec 25 08:39:48 2022\"\n"}
